$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Epoch Accuracy (column B) values per diff
$ws.Range("B5").Value = 0.359375
$ws.Range("B6").Value = 0.296875
$ws.Range("B7").Value = 0.328125
$ws.Range("B8").Value = 0.296875
$ws.Range("B9").Value = 0.28125
$ws.Range("B10").Value = 0.28125
$ws.Range("B11").Value = 0.25
$ws.Range("B12").Value = 0.359375
$ws.Range("B13").Value = 0.328125
$ws.Range("B14").Value = 0.359375
$ws.Range("B15").Value = 0.296875
$ws.Range("B16").Value = 0.1875
$ws.Range("B18").Value = 0.28125
$ws.Range("B21").Value = 0.28125
$ws.Range("B22").Value = 0.234375
$ws.Range("B23").Value = 0.265625
$ws.Range("B24").Value = 0.234375
$ws.Range("B25").Value = 0.265625
$ws.Range("B26").Value = 0.203125
$ws.Range("B28").Value = 0.28125
$ws.Range("B29").Value = 0.234375
$ws.Range("B30").Value = 0.203125
$ws.Range("B31").Value = 0.25
$ws.Range("B32").Value = 0.265625
$ws.Range("B33").Value = 0.328125
$ws.Range("B35").Value = 0.234375
$ws.Range("B36").Value = 0.171875
$ws.Range("B37").Value = 0.171875
$ws.Range("B38").Value = 0.1875
$ws.Range("B39").Value = 0.1875
$ws.Range("B40").Value = 0.1875
$ws.Range("B41").Value = 0.203125
$ws.Range("B42").Value = 0.203125
$ws.Range("B43").Value = 0.203125
$ws.Range("B44").Value = 0.203125
$ws.Range("B45").Value = 0.203125
$ws.Range("B46").Value = 0.1875
$ws.Range("B47").Value = 0.1875
$ws.Range("B48").Value = 0.203125
$ws.Range("B49").Value = 0.203125
$ws.Range("B50").Value = 0.203125
$ws.Range("B62").Value = 0.1875
$ws.Range("B63").Value = 0.1875
$ws.Range("B64").Value = 0.1875
$ws.Range("B65").Value = 0.1875
$ws.Range("B103").Value = 0.125
$ws.Range("B104").Value = 0.140625
$ws.Range("B105").Value = 0.25
$ws.Range("B107").Value = 0.3125
$ws.Range("B108").Value = 0.21875
$ws.Range("B109").Value = 0.1875
$ws.Range("B110").Value = 0.203125
$ws.Range("B112").Value = 0.15625
$ws.Range("B115").Value = 0.1875
$ws.Range("B116").Value = 0.15625
$ws.Range("B117").Value = 0.21875
$ws.Range("B118").Value = 0.2459016393442623

# Update DisplayOutputs object repr text (memory address changed) in column A, rows 102-117
$newAddr = "<__main__.DisplayOutputs object at 0x7f69d02eefd0>"
for ($r = 102; $r -le 117; $r++) {
    $ws.Cells.Item($r, 1).Value = $newAddr
}

